$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 4183.997964392501
$ws.Range("B2").Value = 339.8770259539482
$ws.Range("C2").Value = 145.6179895915934
$ws.Range("D2").Value = 171.6390628165296
$ws.Range("A3").Value = 4802.930819239951
$ws.Range("B3").Value = 284.4180942577647
$ws.Range("C3").Value = 200.3209279978428
$ws.Range("D3").Value = 149.3553366990938
$ws.Range("A4").Value = 5307.628313950999
$ws.Range("B4").Value = 289.8628745767161
$ws.Range("C4").Value = 240.3788832332299
$ws.Range("D4").Value = 153.6825600763833
$ws.Range("A5").Value = 4235.895220995745
$ws.Range("B5").Value = 268.5462458511672
$ws.Range("C5").Value = 193.2232769323498
$ws.Range("D5").Value = 157.9546635178883
$ws.Range("A6").Value = 3207.615118787027
$ws.Range("B6").Value = 217.8627165167376
$ws.Range("C6").Value = 181.7999153081237
$ws.Range("D6").Value = 163.3928629556814
$ws.Range("A7").Value = 5460.01416975
$ws.Range("B7").Value = 275.6419452026545
$ws.Range("C7").Value = 211.2786784888538
$ws.Range("D7").Value = 127.7123672096911
$ws.Range("A8").Value = 4729.444924191721
$ws.Range("B8").Value = 225.636735955958
$ws.Range("C8").Value = 230.7678827308446
$ws.Range("D8").Value = 137.8461951307728
$ws.Range("A9").Value = 5763.487540072081
$ws.Range("B9").Value = 281.5791864892651
$ws.Range("C9").Value = 238.1559315397983
$ws.Range("D9").Value = 132.7885586095901
$ws.Range("A10").Value = 4005.100389293807
$ws.Range("B10").Value = 269.7257029403329
$ws.Range("C10").Value = 177.5508572971512
$ws.Range("D10").Value = 158.9240956328148
$ws.Range("A11").Value = 5316.005401645497
$ws.Range("B11").Value = 347.9975358683374
$ws.Range("C11").Value = 188.6926432053291
$ws.Range("D11").Value = 157.0273936359837
$ws.Range("A12").Value = 4741.710343547656
$ws.Range("B12").Value = 275.440922182457
$ws.Range("C12").Value = 175.1410711063918
$ws.Range("D12").Value = 135.1167439892275
$ws.Range("A13").Value = 4684.33647121434
$ws.Range("B13").Value = 324.3637047326682
$ws.Range("C13").Value = 181.8230192229595
$ws.Range("D13").Value = 164.0814822042272
$ws.Range("A14").Value = 7653.82336983379
$ws.Range("B14").Value = 475.801821374818
$ws.Range("C14").Value = 210.5076500367769
$ws.Range("D14").Value = 148.6673469754573
$ws.Range("A15").Value = 4841.187744523743
$ws.Range("B15").Value = 320.4344663178875
$ws.Range("C15").Value = 134.2856201157396
$ws.Range("D15").Value = 134.0863446006832
$ws.Range("A16").Value = 5488.24291847207
$ws.Range("B16").Value = 325.6988024978548
$ws.Range("C16").Value = 174.8624074659645
$ws.Range("D16").Value = 133.5686808483513
$ws.Range("A17").Value = 5791.573837215479
$ws.Range("B17").Value = 293.866338657123
$ws.Range("C17").Value = 193.9165594124365
$ws.Range("D17").Value = 116.7054470581507
$ws.Range("A18").Value = 6393.957440524506
$ws.Range("B18").Value = 344.8456726550759
$ws.Range("C18").Value = 201.3862045393949
$ws.Range("D18").Value = 124.2906234230692
$ws.Range("A19").Value = 3844.476381341521
$ws.Range("B19").Value = 229.856340720627
$ws.Range("C19").Value = 145.9942001672023
$ws.Range("D19").Value = 130.0867694754717
$ws.Range("A20").Value = 3728.004178736171
$ws.Range("B20").Value = 294.2925835101911
$ws.Range("C20").Value = 188.815513261219
$ws.Range("D20").Value = 185.9080719088816
$ws.Range("A21").Value = 4653.842913897428
$ws.Range("B21").Value = 349.8714702172574
$ws.Range("C21").Value = 149.1868092937349
$ws.Range("D21").Value = 161.9862965404321
$ws.Range("A22").Value = 6159.33942011141
$ws.Range("B22").Value = 336.07380231384
$ws.Range("C22").Value = 197.484190930842
$ws.Range("D22").Value = 126.2475001992578
$ws.Range("A23").Value = 4480.178564410821
$ws.Range("B23").Value = 355.1868118412111
$ws.Range("C23").Value = 115.8455700062837
$ws.Range("D23").Value = 154.6888663777746
$ws.Range("A24").Value = 6024.831251206691
$ws.Range("B24").Value = 301.704388040833
$ws.Range("C24").Value = 254.7585125977538
$ws.Range("D24").Value = 141.5114800718186
$ws.Range("A25").Value = 4996.398591338498
$ws.Range("B25").Value = 248.6627862687093
$ws.Range("C25").Value = 204.1283517919147
$ws.Range("D25").Value = 127.1776018700035
$ws.Range("A26").Value = 5488.700932565435
$ws.Range("B26").Value = 293.27182129078
$ws.Range("C26").Value = 198.2391720952879
$ws.Range("D26").Value = 128.9359866278049
$ws.Range("A27").Value = 3874.91338159789
$ws.Range("B27").Value = 219.5366126256114
$ws.Range("C27").Value = 197.7313522844315
$ws.Range("D27").Value = 148.7248245086496
$ws.Range("A28").Value = 3693.457543500582
$ws.Range("B28").Value = 250.2610594420158
$ws.Range("C28").Value = 185.3177663406039
$ws.Range("D28").Value = 164.0272587594507
$ws.Range("A29").Value = 5934.265389824087
$ws.Range("B29").Value = 346.4583244888966
$ws.Range("C29").Value = 219.9231786735873
$ws.Range("D29").Value = 149.8274646099043
$ws.Range("A30").Value = 4193.924967920159
$ws.Range("B30").Value = 187.3633410665345
$ws.Range("C30").Value = 237.4118350135104
$ws.Range("D30").Value = 140.9670930739501
$ws.Range("A31").Value = 6191.927984488466
$ws.Range("B31").Value = 326.5676347170769
$ws.Range("C31").Value = 152.8317395536022
$ws.Range("D31").Value = 99.25367650454592
$ws.Range("A32").Value = 4033.369442941071
$ws.Range("B32").Value = 308.1219828907444
$ws.Range("C32").Value = 173.4496534425595
$ws.Range("D32").Value = 174.7155044466616
$ws.Range("A33").Value = 5076.778813891905
$ws.Range("B33").Value = 235.0118107217378
$ws.Range("C33").Value = 179.5008117926838
$ws.Range("D33").Value = 105.9910648012135
$ws.Range("A34").Value = 4876.407776977804
$ws.Range("B34").Value = 270.7311045730235
$ws.Range("C34").Value = 215.3827214584132
$ws.Range("D34").Value = 147.3852631070871
$ws.Range("A35").Value = 4475.125806226122
$ws.Range("B35").Value = 253.0503061067726
$ws.Range("C35").Value = 189.4515496103742
$ws.Range("D35").Value = 140.3053789367162
$ws.Range("A36").Value = 6145.494816287224
$ws.Range("B36").Value = 322.9751339606493
$ws.Range("C36").Value = 242.3866251370197
$ws.Range("D36").Value = 141.7593500431873
$ws.Range("A37").Value = 5604.988281942416
$ws.Range("B37").Value = 362.7132194412006
$ws.Range("C37").Value = 210.673783076294
$ws.Range("D37").Value = 164.6694187592871
$ws.Range("A38").Value = 4718.203433362023
$ws.Range("B38").Value = 242.8393186615144
$ws.Range("C38").Value = 197.6363103605439
$ws.Range("D38").Value = 130.8650743850901
$ws.Range("A39").Value = 3814.495081269925
$ws.Range("B39").Value = 196.1254644536162
$ws.Range("C39").Value = 233.7527555131854
$ws.Range("D39").Value = 156.5888064997174
$ws.Range("A40").Value = 5723.613259955642
$ws.Range("B40").Value = 358.2460482821252
$ws.Range("C40").Value = 213.0661590240718
$ws.Range("D40").Value = 159.4914480582309
$ws.Range("A41").Value = 5426.333359442152
$ws.Range("B41").Value = 287.9275050006977
$ws.Range("C41").Value = 198.3184561001557
$ws.Range("D41").Value = 128.6231503496672
$ws.Range("A42").Value = 6854.416276612737
$ws.Range("B42").Value = 339.5087528826947
$ws.Range("C42").Value = 240.5488235753697
$ws.Range("D42").Value = 124.3829819925229
$ws.Range("A43").Value = 5519.908893505078
$ws.Range("B43").Value = 284.2602669551711
$ws.Range("C43").Value = 195.3722818791551
$ws.Range("D43").Value = 122.1947838367589
$ws.Range("A44").Value = 5812.679976960413
$ws.Range("B44").Value = 393.0155345617394
$ws.Range("C44").Value = 176.5530659454782
$ws.Range("D44").Value = 155.9591666594987
$ws.Range("A45").Value = 2358.791032646498
$ws.Range("B45").Value = 199.5836013900623
$ws.Range("C45").Value = 166.2960821535796
$ws.Range("D45").Value = 176.4190290882778
$ws.Range("A46").Value = 2805.316522659203
$ws.Range("B46").Value = 210.294318940217
$ws.Range("C46").Value = 168.4182881983935
$ws.Range("D46").Value = 167.2836082406685
$ws.Range("A47").Value = 4799.66348757427
$ws.Range("B47").Value = 327.8877885792913
$ws.Range("C47").Value = 176.0657656101837
$ws.Range("D47").Value = 159.113266666036
$ws.Range("A48").Value = 5481.531930333672
$ws.Range("B48").Value = 275.2080967925227
$ws.Range("C48").Value = 190.5025549415579
$ws.Range("D48").Value = 116.6524120156672
$ws.Range("A49").Value = 4446.618750085127
$ws.Range("B49").Value = 273.760029201124
$ws.Range("C49").Value = 161.7423283605575
$ws.Range("D49").Value = 138.2471600005978
$ws.Range("A50").Value = 6067.398239699472
$ws.Range("B50").Value = 300.8542034914324
$ws.Range("C50").Value = 211.1937319486963
$ws.Range("D50").Value = 118.8414526414741
$ws.Range("A51").Value = 6322.9830926214
$ws.Range("B51").Value = 331.2980776498941
$ws.Range("C51").Value = 232.2640057816323
$ws.Range("D51").Value = 134.6881246898658
$ws.Range("A52").Value = 4603.535073643502
$ws.Range("B52").Value = 305.6844418426081
$ws.Range("C52").Value = 188.5992199126504
$ws.Range("D52").Value = 161.056441284989
$ws.Range("A53").Value = 4462.315833996377
$ws.Range("B53").Value = 284.9842960683523
$ws.Range("C53").Value = 166.2319996010016
$ws.Range("D53").Value = 145.2716515661839
$ws.Range("A54").Value = 5303.495755183319
$ws.Range("B54").Value = 336.9581147944488
$ws.Range("C54").Value = 197.7166555421311
$ws.Range("D54").Value = 156.4058941601969
$ws.Range("A55").Value = 4317.890379467585
$ws.Range("B55").Value = 335.5516507371507
$ws.Range("C55").Value = 150.3272692066299
$ws.Range("D55").Value = 167.0918534003308
$ws.Range("A56").Value = 4825.072007925428
$ws.Range("B56").Value = 288.4750339675315
$ws.Range("C56").Value = 156.8253732076046
$ws.Range("D56").Value = 129.9134243263363
$ws.Range("A57").Value = 5535.571688693568
$ws.Range("B57").Value = 312.9150426850317
$ws.Range("C57").Value = 156.4252694920661
$ws.Range("D57").Value = 117.0639131423928
$ws.Range("A58").Value = 4350.507777149798
$ws.Range("B58").Value = 242.0820239588815
$ws.Range("C58").Value = 225.9469376244147
$ws.Range("D58").Value = 156.6731375042553
$ws.Range("A59").Value = 4760.352965150912
$ws.Range("B59").Value = 324.3416426696068
$ws.Range("C59").Value = 215.7404804041033
$ws.Range("D59").Value = 177.5775964634024
$ws.Range("A60").Value = 6387.986441684528
$ws.Range("B60").Value = 389.3133813468136
$ws.Range("C60").Value = 209.2575271834944
$ws.Range("D60").Value = 149.9533731874854
$ws.Range("A61").Value = 4984.874357122745
$ws.Range("B61").Value = 340.1113446044977
$ws.Range("C61").Value = 185.9353666483861
$ws.Range("D61").Value = 163.2308101898517
